$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 5 de Octubre de 2020 a las 05:11"

# 2. Row 35 - Belgica: refresh case counts
$ws.Range("B35").Value = 130235
$ws.Range("C35").Value = 2612
$ws.Range("D35").Value = 19679
$ws.Range("E35").Value = 100492
$ws.Range("G35").Value = 20
$ws.Range("H35").Value = 10064

# 3. Row 39 - Kazajistan: refresh case counts
$ws.Range("B39").Value = 108236
$ws.Range("C39").Value = 59
$ws.Range("D39").Value = 103277
$ws.Range("E39").Value = 3234

# 4. Row 81 - Australia: refresh case counts
$ws.Range("B81").Value = 27144
$ws.Range("C81").Value = 8
$ws.Range("E81").Value = 1384

# 5. Row 133 - Trinidad yTobago: refresh case counts
$ws.Range("B133").Value = 4763
$ws.Range("E133").Value = 1839

# 6. Rows 215/216 - reorder "Islas Malvinas" ahead of "Montserrat"
#    (swap the two full rows, matching the new shared-string ordering)
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 13
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0

$ws.Range("A216").Value = "Montserrat"
$ws.Range("B216").Value = 13
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 12
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 1
